$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 'What is the name of the first animal to land on the moon?'
$ws.Range("D2").Value = 'No animal has ever landed on the moon yet.'
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = '{''Explanation:'': ''No animals, including dogs, were sent to the moon. However, if we consider humans as animals, then Neil Armstrong was the first animal (human) to land on the moon.'', ''Answer:'': ''Neil Armstrong (human being, part of animal kingdom)'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.>'', ''Core Concept:'': ''Animal on the moon'', ''Premise of the Question:'': ''Invalid. The question assumes that a non-human animal was sent to the moon, which is not true.''}'
$ws.Range("G2").Value = 'It seems the candidate responses could not reach an agreement for self-consistency to work.'
$ws.Range("H2").Value = -1
$ws.Range("I2").Value = '{''QueryID'': ''0'', ''answer_box'': [{''name'': ''What was the first animal to land on the moon?'', ''url'': ''https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.'', ''snippet'': ''Answer and Explanation: No animals were ever sent to the Moon. Although, since humans are technically animals, one could say that the first animal sent to the Moon was Neil Armstrong. He belonged to the species Homo sapiens.''}], ''related_questions'': [{''name'': ''What was the first animal to come on land?'', ''url'': ''https://bigthink.com/hard-science/first-creature-to-live-on-land/#:~:text=Scientists%20think%20an%20insect%20similar,first%2Dever%20land%2Ddweller.&text=An%20ancient%20millipede%2Dlike%20creature,creature%20to%20live%20on%20land.'', ''snippet'': ''Scientists think an insect similar to the modern millipede crawled around Scotland 425 million years ago, making it the first-ever land-dweller. An ancient millipede-like creature living in Scotland may have been the first creature to live on land.''}, {''name'': ''What was the name of the first land on the moon?'', ''url'': "https://www.rmg.co.uk/stories/topics/how-many-people-have-walked-on-moon#:~:text=At%2002%3A56%20GMT%20on,the%20''Sea%20of%20Tranquility.''", ''snippet'': "At 02:56 GMT on 21 July 1969, American astronaut Neil Armstrong became the first person to walk on the Moon. He stepped out of the Apollo 11 lunar module and onto the Moon''s surface, in an area called the ''Sea of Tranquility.''"}, {''name'': ''What were the first animals in space?'', ''url'': ''https://www.space.com/animals-in-space#:~:text=The%20first%20animals%20to%20reach,flights%20for%20about%20a%20decade.'', ''snippet'': ''The first animals to reach space were fruit flies that the United States launched aboard captured German rockets in 1947. The first mammal to reach space was a rhesus monkey named Albert II, who flew two years later. Both these missions were suborbital, as were all animal flights for about a decade.''}, {''name'': ''Was a dog the first on the moon?'', ''url'': ''https://www.quora.com/Why-did-they-send-a-dog-to-the-Moon-first-Why-didnt-they-choose-another-animal#:~:text=No%20dog%20has%20ever%20been,on%20the%20Soviet%20Sputnik%202.'', ''snippet'': ''No dog has ever been sent to the moon. The first (deliberate*) living thing launched into orbit was a dog name Laika, on the Soviet Sputnik 2.''}], ''organic_results'': [{''name'': ''What was the first animal in space?'', ''url'': ''https://www.rmg.co.uk/stories/topics/what-was-first-animal-space'', ''snippet'': ''The first animals to orbit the moon and return to Earth were two Russian tortoises aboard Zond 5. On 15 September 1968, the tortoises were launched with plants, ...''}, {''name'': ''Is it true that the first animal that went to the moon was a dog?'', ''url'': ''https://www.quora.com/Is-it-true-that-the-first-animal-that-went-to-the-moon-was-a-dog'', ''snippet'': ''No dog has ever been sent to the moon. The first (deliberate*) living thing launched into orbit was a dog name Laika, on the Soviet Sputnik 2.''}]}'
$ws.Range("J2").Value = '{''Explanation:'': ''No animals, including dogs, were sent to the moon. However, if we consider humans as animals, then Neil Armstrong was the first animal (homo sapien) to land on the moon.'', ''Answer:'': ''Neil Armstrong, considering humans as animals.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens>'', ''Core Concept:'': ''First animal on the moon'', ''Premise of the Question:'': ''Invalid. The question is based on a false premise that an animal landed on the moon.'', ''Certainty_Estimation'': ''Yes.\n\nThe two responses are discussing the same topic. Both Response A and Response B mention Neil Armstrong and the fact that he was the first being to land on the moon, considering humans as animals (homo sapiens). The context and key points are the same in both responses, with the only difference being the specific term used for humans in Response B.''}'
$ws.Range("K2").Value = '{''Explanation:'': ''No animals, including dogs, were ever sent to the moon. However, if we consider humans as animals, one could say that the first animal to reach the moon was Neil Armstrong.'', ''Answer:'': ''Neil Armstrong (considering humans as animals)'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.>'', ''Core Concept:'': ''Space exploration'', ''Premise of the Question:'': ''Invalid. The premise assumes that an animal other than humans landed on the moon, but this is not true.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("L2").Value = '{''Explanation:'': ''No animals, including dogs, were ever sent to the moon. However, the first animal sent to space was a fruit fly.'', ''Answer:'': ''No animals have landed on the moon.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://www.space.com/animals-in-space>'', ''Core Concept:'': ''Space Animals'', ''Premise of the Question:'': ''Invalid. The question contains a false premise that an animal landed on the moon.'', ''Certainty_Estimation'': ''No.\n\nConfidence: 90%''}'
$ws.Range("M2").Value = '{''Explanation:'': ''No animals, including dogs, were sent to the moon. The first living beings in space were fruit flies and the first mammal in space was a rhesus monkey named Albert II.'', ''Answer:'': ''No animals have landed on the moon.'', ''Confidence Level:'': ''90%'', ''Source:'': ''Multiple sources from the external links provided.'', ''Core Concept:'': ''Animal in space.'', ''Premise of the Question:'': ''Invalid. The question contains a false premise that an animal landed on the moon.'', ''Certainty_Estimation'': ''No.\n\nThe two responses discuss different aspects of space exploration and the presence of animals in space. Response A focuses on the fact that Neil Armstrong was the first human (and thus, an animal) to land on the moon, while Response B clarifies that no animals have landed on the moon and mentions the first animals sent to space.''}'
$ws.Range("N2").Value = '{''Explanation:'': ''No animals, including dogs, were ever sent to the Moon. However, there were animals sent to space before humans.'', ''Answer:'': ''No animals have landed on the Moon.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens>'', ''Core Concept:'': ''Animal, Moon'', ''Premise of the Question:'': ''Invalid. The information provided in the question is not accurate. Dogs and other animals have not landed on the Moon.'', ''Certainty_Estimation'': ''No.\n\nThe two responses discuss related topics but are not discussing the same topic. Response A focuses on defining Neil Armstrong as the first "animal" (human) on the moon, while Response B clarifies that no animals, including humans, have landed on the moon, but some animals were sent to space.''}'
$ws.Range("O2").Value = 'Incorrect'
$ws.Range("P2").Value = 'The response does not provide any information or answer to the question.'

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 'What is the name of Leonardo DiCaprio''s third child?'
$ws.Range("D3").Value = 'Leonardo DiCaprio does not have any children.'
$ws.Range("E3").Value = 'Leonardo DiCaprio does not have any children. The question asks about Leonardo DiCaprio''s third child, but based on the information from reliable sources, it is confirmed that DiCaprio does not have any children.'
$ws.Range("F3").Value = '{''Explanation:'': "The question asks about Leonardo DiCaprio''s third child, but based on the information from reliable sources, DiCaprio does not have any children.", ''Answer:'': ''Leonardo DiCaprio does not have any children.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have>'', ''Core Concept:'': "Leonardo DiCaprio''s children", ''Premise of the Question:'': ''Invalid. The information provided in the question is not accurate as Leonardo DiCaprio does not have any children.''}'
$ws.Range("G3").Value = '{''Explanation:'': "The question asks about Leonardo DiCaprio''s third child, but based on the information from reliable sources, it is confirmed that DiCaprio does not have any children.", ''Answer:'': ''Leonardo DiCaprio does not have any children.'', ''Confidence Level:'': ''100.0%'', ''Source:'': ''Quora, Wikipedia, Yahoo Entertainment, kahawatungu'', ''Core Concept:'': "Leonardo DiCaprio''s children", ''Premise of the Question:'': ''Invalid. The information provided in the question is incorrect as Leonardo DiCaprio does not have any children.'', ''Certainty_Estimation'': "Yes, both responses are discussing the same topic. They both address the question about Leonardo DiCaprio''s third child and confirm that he does not have any children."}'
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = '{''QueryID'': ''1'', ''related_questions'': [{''name'': ''How many kids does Leo DiCaprio have?'', ''url'': ''https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have'', ''snippet'': ''He has never been married and has no children. How many relationships did Leonardo DiCaprio had in the past? Bridget Hall (1994): The model was the first woman linked to DiCaprio after the start of his movie career. Kristen Zang (1996): DiCaprio and Zang spent time together before his meteoric rise with Titanic.''}, {''name'': "What is Leonardo DiCaprio''s birth name?", ''url'': ''https://en.wikipedia.org/wiki/Leonardo_DiCaprio'', ''snippet'': ''Leonardo Wilhelm DiCaprio (/diˈkæprioʊ, dɪ-/; Italian: [diˈkaːprjo]; born November 11, 1974) is an American actor and film producer.''}, {''name'': "What is Leonardo DiCaprio''s longest relationship?", ''url'': ''https://www.yahoo.com/entertainment/leonardo-dicaprio-exes-said-dating-210507180.html'', ''snippet'': "One of Leonardo DiCaprio''s longest relationships was with Gisele Bundchen, whom he dated from 2001 to 2005. Bundchen has spoken about their relationship the most out of all of his exes. In a 2009 interview with Vanity Fair, she spoke about it, and had quite nice things to say."}, {''name'': ''How many siblings does Leonardo DiCaprio have?'', ''url'': ''https://kahawatungu.com/leonardo-dicaprio-siblings/'', ''snippet'': "DiCaprio''s parents divorced when he was a year old, and he was raised by his mother. He does not have any biological siblings, but he has a stepbrother named Adam Farrar. The relationship between DiCaprio and Farrar has been strained, and they have reportedly not been in contact for years."}], ''organic_results'': [{''name'': ''Leonardo DiCaprio'', ''url'': ''https://en.wikipedia.org/wiki/Leonardo_DiCaprio'', ''snippet'': "Leonardo Wilhelm DiCaprio was born on November 11, 1974, in Los Angeles, California. · DiCaprio''s parents named him Leonardo because his pregnant mother first ..."}, {''name'': "What Kate Winslet''s Children Call Leonardo DiCaprio", ''url'': ''https://www.nickiswift.com/455252/what-kate-winslets-children-call-leonardo-dicaprio/'', ''snippet'': "Leonardo DiCaprio and Kate Winslet''s unsinkable friendship ... These days, Kate Winslet is a proud mama of three children. She welcomed Mia Honey ..."}]}'
$ws.Range("J3").Value = '{''Explanation:'': "The question asks for the name of Leonardo DiCaprio''s third child. However, based on the provided sources, Leonardo DiCaprio does not have any children.", ''Answer:'': ''Leonardo DiCaprio does not have a third child because he does not have any children.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<Quora, Wikipedia, Yahoo, kahawatungu>'', ''Core Concept:'': ''Children of Leonardo DiCaprio'', ''Premise of the Question:'': ''Invalid. The question contains a false premise as Leonardo DiCaprio does not have any children.'', ''Certainty_Estimation'': "Yes.\n\nBoth responses are discussing the same topic, which is Leonardo DiCaprio''s children or lack thereof. The main subject, context, and key points in both responses are consistent with each other, emphasizing that DiCaprio does not have any children, and therefore, he does not have a third child."}'
$ws.Range("K3").Value = '{''Explanation:'': "The question asks about Leonardo DiCaprio''s third child, but based on the provided sources, there is no information that DiCaprio has any children.", ''Answer:'': ''Leonardo DiCaprio does not have a third child because he does not have any children at all.'', ''Confidence Level:'': ''90%'', ''Source:'': ''Multiple sources including Wikipedia, Quora, Yahoo Entertainment, and Kahawa Tungu (<https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have>, <https://en.wikipedia.org/wiki/Leonardo_DiCaprio>, <https://www.yahoo.com/entertainment/leonardo-dicaprio-exes-said-dating-210507180.html>, <https://kahawatungu.com/leonardo-dicaprio-siblings/>)'', ''Core Concept:'': "Leonardo DiCaprio''s children", ''Premise of the Question:'': ''Invalid. The premise that DiCaprio has a third child is incorrect based on the available information.'', ''Certainty_Estimation'': ''Yes.\n\nConfidence: 100%''}'
$ws.Range("L3").Value = '{''Explanation:'': "The question asks about Leonardo DiCaprio''s third child, but based on the information from reliable sources, it is confirmed that DiCaprio does not have any children.", ''Answer:'': ''Leonardo DiCaprio does not have any children.'', ''Confidence Level:'': ''95%'', ''Source:'': ''Quora, Wikipedia, Yahoo Entertainment, kahawatungu'', ''Core Concept:'': "Leonardo DiCaprio''s children", ''Premise of the Question:'': ''Invalid. The information provided in the question is incorrect as Leonardo DiCaprio does not have any children.'', ''Certainty_Estimation'': "Yes, both responses are discussing the same topic. They both address the question about Leonardo DiCaprio''s third child and confirm that he does not have any children."}'
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").Value = 'Correct'
$ws.Range("P3").Value = 'The response accurately addresses the question by stating that Leonardo DiCaprio does not have any children, which directly answers the query about his third child. Therefore, the response is correct.'

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 'What year did the first human land on Mars?'
$ws.Range("D4").Value = 'No humans have been to Mars yet.'
$ws.Range("E4").Value = 'Humans have not yet landed on Mars. The question asks for the year of the first human landing on Mars, but as of now, humans have not landed on Mars. The first successful Mars landing was by the Soviet Union''s Mars 3 in 1971, followed by successful landings by the United States and China. However, these were not human landings.'
$ws.Range("F4").Value = '{''Explanation:'': "The question asks for the year a human first landed on Mars, but there is no evidence to suggest that a human has ever landed on Mars. The first successful spacecraft landing on Mars was in 1971 by the Soviet Union''s Mars 3.", ''Answer:'': ''No human has landed on Mars yet.'', ''Confidence Level:'': ''90%'', ''Source:'': "<https://en.wikipedia.org/wiki/Mars_landing#:~:text=Soviet%20Union''s%20Mars%203%2C%20which,have%20conducted%20Mars%20landings%20successfully.>", ''Core Concept:'': ''Mars landing'', ''Premise of the Question:'': ''Invalid. The question contains a false premise that a human has already landed on Mars.''}'
$ws.Range("G4").Value = '{''Explanation:'': "The question asks for the year of the first human landing on Mars, but as of now, humans have not landed on Mars. The first successful Mars landing was by the Soviet Union''s Mars 3 in 1971, followed by successful landings by the United States and China. However, these were not human landings.", ''Answer:'': ''Humans have not yet landed on Mars.'', ''Confidence Level:'': ''100.0%'', ''Source:'': ''<My knowledge>'', ''Core Concept:'': ''Mars landing'', ''Premise of the Question:'': ''Invalid. The question contains a false premise as humans have not landed on Mars yet.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = '{''QueryID'': ''2'', ''related_questions'': [{''name'': ''When did man first land on Mars?'', ''url'': "https://en.wikipedia.org/wiki/Mars_landing#:~:text=Soviet%20Union''s%20Mars%203%2C%20which,have%20conducted%20Mars%20landings%20successfully.", ''snippet'': "Soviet Union''s Mars 3, which landed in 1971, was the first successful Mars landing. As of 2023, the Soviet Union, United States and China have conducted Mars landings successfully."}, {''name'': ''What year will the first person be on Mars?'', ''url'': ''https://en.wikipedia.org/wiki/Human_mission_to_Mars#:~:text=The%20first%20crewed%20Mars%20Mission,is%20proposed%20for%20the%202030s.'', ''snippet'': ''The first crewed Mars Mission, which would include sending astronauts to Mars, orbiting Mars, and a return to Earth, is proposed for the 2030s.''}, {''name'': ''Who went 1st on Mars?'', ''url'': "https://en.wikipedia.org/wiki/Exploration_of_Mars#:~:text=Starting%20in%201960%2C%20the%20Soviets,1965%2C%20by%20NASA''s%20Mariner%204.", ''snippet'': "Starting in 1960, the Soviets launched a series of probes to Mars including the first intended flybys and hard (impact) landing (Mars 1962B). The first successful flyby of Mars was on 14–15 July 1965, by NASA''s Mariner 4."}, {''name'': ''When was the first US Mars landing?'', ''url'': ''https://www.planetary.org/space-missions/viking#:~:text=Viking%201%20launched%20on%20August,touched%20down%20on%20September%203.'', ''snippet'': ''Viking 1 launched on August 20, 1975. It arrived in Mars orbit on June 19, 1976 and the lander touched down on July 20, 1976. Viking 2 launched less than a month after Viking 1 on September 9, 1975. The spacecraft arrived in orbit on August 7, 1976 and the lander touched down on September 3.''}], ''organic_results'': [{''name'': ''Human mission to Mars'', ''url'': ''https://en.wikipedia.org/wiki/Human_mission_to_Mars'', ''snippet'': ''The first crewed Mars Mission, which would include sending astronauts to Mars, orbiting Mars, and a return to Earth, is proposed for the 2030s.''}, {''name'': ''Mars Exploration Timeline - the NSSDCA'', ''url'': ''https://nssdc.gsfc.nasa.gov/planetary/chronology_mars.html'', ''snippet'': ''Mission Timeline ; Mariner 6 - 25 February 1969 - Mars Flyby ; Mariner 7 - 27 March 1969 - Mars Flyby''}]}'
$ws.Range("J4").Value = '{''Explanation:'': "The question asks for the year a human first landed on Mars, but according to the provided sources, no human has landed on Mars yet. The first successful Mars landing was by the Soviet Union''s Mars 3 in 1971, and the first US Mars landing was by the Viking 1 lander in 1976, but these are not human landings.", ''Answer:'': ''No human has landed on Mars yet.'', ''Confidence Level:'': ''100%'', ''Source:'': "<https://en.wikipedia.org/wiki/Mars_landing#:~:text=Soviet%20Union''s%20Mars%203%2C%20which,have%20conducted%20Mars%20landings%20successfully.>", ''Core Concept:'': ''Mars landing'', ''Premise of the Question:'': ''Invalid. The question contains a false premise, as no human has landed on Mars yet.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("K4").Value = 'It seems the candidate response 1 was missing some keys in the response dict {''message'': "Explanation: According to the information provided, no human has landed on Mars yet. The first successful flyby of Mars was in 1965 by NASA''s Mariner 4, and the first US Mars landing was in 1976 with the Viking 1 mission. However, there are plans for a crewed Mars mission in the 2"} so the current                       iteration of the candidate response generation has been skipped. The next iteration                       will continue.'
$ws.Range("L4").Value = '{''Explanation:'': "The question asks for the year of the first human landing on Mars, but as of now, humans have not landed on Mars. The first successful Mars landing was by the Soviet Union''s Mars 3 in 1971, followed by successful landings by the United States and China. However, these were not human landings.", ''Answer:'': ''Humans have not yet landed on Mars.'', ''Confidence Level:'': ''100%'', ''Source:'': ''<My knowledge>'', ''Core Concept:'': ''Mars landing'', ''Premise of the Question:'': ''Invalid. The question contains a false premise as humans have not landed on Mars yet.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("M4").Value = '{''Explanation:'': "The question asks for the first human landing on Mars, and according to the provided sources, the first successful Mars landing was conducted by the Soviet Union''s Mars 3 in 1971. However, this was an unmanned mission, and no human was on board. Therefore, the answer to the question is that no human has landed on Mars yet.", ''Answer:'': ''No human has landed on Mars yet.'', ''Confidence Level:'': ''90%'', ''Source:'': ''{"QueryID": "2", "related_questions": [...], "organic_results": [...]}, <https://en.wikipedia.org/wiki/Mars_landing#:~:text=Soviet%20Union\''s%20Mars%203%2C%20which,have%20conducted%20Mars%20landings%20successfully.>'', ''Core Concept:'': ''Mars landing'', ''Premise of the Question:'': ''Valid, but the information provided in the question is incorrect.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("N4").ClearContents()
$ws.Range("O4").Value = 'Correct'
$ws.Range("P4").Value = 'The response accurately states that no humans have landed on Mars yet, which is the correct answer to the question. The additional information about Mars landings does not contradict the primary answer.'

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 'What was the name of the Zodiac killer?'
$ws.Range("D5").Value = 'The Zodiac Killer has not been identified.'
$ws.Range("E5").Value = 'The most famous suspect is Arthur Leigh Allen, but his identity was never confirmed. The Zodiac Killer''s true identity remains unknown, but Arthur Leigh Allen is the most well-known suspect. Although there is significant circumstantial evidence against him, Allen was never officially charged, and his DNA does not match any found at the crime scenes.'
$ws.Range("F5").Value = '{''Explanation:'': "The Zodiac Killer''s true identity remains unknown, despite investigations and theories. Arthur Leigh Allen is one of the most well-known suspects, but he was never definitively proven to be the killer.", ''Answer:'': "The Zodiac Killer''s identity is unknown, and Arthur Leigh Allen is one of the main suspects.", ''Confidence Level:'': ''85%'', ''Source:'': ''Multiple sources, including <https://www.biography.com/crime/zodiac-killer> and <https://en.wikipedia.org/wiki/Zodiac_Killer>'', ''Core Concept:'': "Zodiac Killer''s identity", ''Premise of the Question:'': "The premise is invalid, as the Zodiac Killer''s identity is not publicly known."}'
$ws.Range("G5").Value = '{''Explanation:'': "The Zodiac Killer''s true identity remains unknown, but Arthur Leigh Allen is the most well-known suspect. Although there is significant circumstantial evidence against him, Allen was never officially charged, and his DNA does not match any found at the crime scenes.", ''Answer:'': ''The most famous suspect is Arthur Leigh Allen, but his identity was never confirmed.'', ''Confidence Level:'': ''100.0%'', ''Source:'': ''<https://www.biography.com/crime/zodiac-killer#:~:text>'', ''Core Concept:'': ''Zodiac Killer identity'', ''Premise of the Question:'': ''Valid, as the Zodiac Killer was a real serial killer, but the premise about his name is invalid because his true identity remains unknown.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = '{''QueryID'': ''3'', ''answer_box'': [{''name'': ''Zodiac Killer: Biography, Serial Killer, Criminal'', ''url'': ''https://www.biography.com/crime/zodiac-killer#:~:text=True%2Dcrime%20author%20and%20former,as%20the%20most%20likely%20suspect.'', ''snippet'': "True-crime author and former San Francisco Chronicle cartoonist Robert Graysmith wrote two separate works on the killer (1986''s Zodiac and 2002''s Zodiac Unmasked), ultimately identifying a man named Arthur Leigh Allen as the most likely suspect."}], ''related_questions'': [{''name'': "What is the Zodiac Killer''s real name?", ''url'': ''https://en.wikipedia.org/wiki/Zodiac_Killer#:~:text=Of%20the%20four%20ciphers%20he,offender%20who%20died%20in%201992.'', ''snippet'': ''Of the four ciphers he produced, two remain unsolved, and one was cracked only in 2020. While many theories regarding the identity of the killer have been suggested, the only suspect authorities ever publicly named was Arthur Leigh Allen, a former elementary school teacher and convicted sex offender who died in 1992.''}, {''name'': ''Is Arthur Leigh Allen the Zodiac?'', ''url'': ''https://screenrant.com/all-evidence-arthur-leigh-allen-not-zodiac-killer/#:~:text=Arthur%20Leigh%20Allen%20was%20the,the%20case%20is%20still%20open.'', ''snippet'': ''Arthur Leigh Allen was the most likely suspect in the Zodiac case, but there was not enough physical evidence to prove his guilt. The circumstantial evidence and fear in the community led to Allen being falsely identified as the killer. The identity of the Zodiac Killer remains a mystery, and the case is still open.''}, {''name'': ''Were there 2 Zodiac killers?'', ''url'': ''https://www.distractify.com/p/how-many-zodiac-killers-were-there#:~:text=When%20Choon%2DJae%20was%20later,anything%20is%20possible.'', ''snippet'': ''When Choon-Jae was later caught in 2018, he could no longer be compared to the one and only Zodiac. So there are at least three other killers who have been compared to the Zodiac Killer, and while experts agree that the original Zodiac Killer is likely just one man … anything is possible.''}, {''name'': ''Did they ever catch Jack the Ripper?'', ''url'': ''https://www.britannica.com/biography/Jack-the-Ripper#:~:text=Jack%20the%20Ripper%20was%20an,macabre%20tourist%20industry%20in%20London.'', ''snippet'': "Jack the Ripper was an English serial killer. Between August and November 1888, he murdered at least five women—all prostitutes—in or near the Whitechapel district of London''s East End. Jack the Ripper was never identified or arrested. Today the murder sites are the locus of a macabre tourist industry in London."}], ''organic_results'': [{''name'': ''Zodiac Killer'', ''url'': ''https://en.wikipedia.org/wiki/Zodiac_Killer'', ''snippet'': ''For the Japanese Zodiac copycat, see Kobe child murders. The Zodiac Killer is the pseudonym of an unidentified serial killer who operated in Northern California ...''}, {''name'': "Zodiac Killer''s identity has been reportedly revealed, DNA ...", ''url'': ''https://www.marca.com/en/lifestyle/us-news/2023/05/20/64680e6de2704e873f8b4591.html'', ''snippet'': "The Zodiac Killer''s identity has been reportedly found by the FBI, they believe that although the infamous murderer that was active during ..."}]}'
$ws.Range("J5").Value = '{''Explanation:'': "The Zodiac Killer''s true identity remains unknown, despite many suspicions and theories, with Arthur Leigh Allen being the most widely-known suspect. Although the FBI recently reportedly identified the Zodiac Killer, the name has not been officially released, leaving the mystery still unsolved.", ''Answer:'': ''The true identity of the Zodiac Killer remains unknown.'', ''Confidence Level:'': ''85%'', ''Source:'': ''<https://en.wikipedia.org/wiki/Zodiac_Killer>, <https://www.marca.com/en/lifestyle/us-news/2023/05/20/64680e6de2704e873f8b4591.html>'', ''Core Concept:'': "Zodiac Killer''s Identity", ''Premise of the Question:'': ''Invalid. The question assumes there is a known name for the Zodiac Killer, but no conclusive evidence has led to a verified identity.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("K5").Value = '{''Explanation:'': "The Zodiac Killer is an unidentified serial killer who operated in Northern California during the late 1960s and early 1970s. Although there have been many theories about the killer''s identity, and one main suspect named Arthur Leigh Allen, the case remains open and the killer''s identity has never been officially confirmed. A recent report claiming that the FBI has identified the Zodiac Killer is not yet officially confirmed.", ''Answer:'': ''The identity of the Zodiac Killer remains officially unconfirmed.'', ''Confidence Level:'': ''85%'', ''Source:'': ''Multiple sources, including FBI, Britannica, and Wikipedia. (<https://www.fbi.gov/history/famous-cases/zodiac-killer>; <https://www.britannica.com/biography/Jack-the-Ripper>; <https://en.wikipedia.org/wiki/Zodiac_Killer>)'', ''Core Concept:'': ''Zodiac Killer'', ''Premise of the Question:'': ''The premise is valid, as it asks for the identity of the Zodiac Killer, but the answer is that the identity remains officially unconfirmed.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("L5").Value = '{''Explanation:'': "The Zodiac Killer''s true identity remains unknown, but Arthur Leigh Allen is the most well-known suspect. Although there is significant circumstantial evidence against him, Allen was never officially charged, and his DNA does not match any found at the crime scenes.", ''Answer:'': ''The most famous suspect is Arthur Leigh Allen, but his identity was never confirmed.'', ''Confidence Level:'': ''85%'', ''Source:'': ''<https://www.biography.com/crime/zodiac-killer#:~:text>'', ''Core Concept:'': ''Zodiac Killer identity'', ''Premise of the Question:'': ''Valid, as the Zodiac Killer was a real serial killer, but the premise about his name is invalid because his true identity remains unknown.'', ''Certainty_Estimation'': ''Yes''}'
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").Value = 'Correct'
$ws.Range("P5").Value = 'The primary answer in the response correctly indicates that the Zodiac Killer''s true identity remains unknown, which aligns with the correct answer. The additional information about Arthur Leigh Allen being a suspect does not contradict the primary answer or reshape one''s perception of it. Therefore, the response is credited.'

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 'Why are all quickly verifiable problems also quickly solvable?'
$ws.Range("D6").Value = 'This is the P versus NP problem and it remains open.'
$ws.Range("E6").Value = 'Not all quickly verifiable problems are also quickly solvable. The premise of the question assumes that if a problem can be quickly verified, then it can also be quickly solved. However, this is the fundamental question that the P versus NP problem tries to answer. In the field of computer science, a problem is considered to be in NP if a potential solution can be verified quickly, but not necessarily solved quickly. Therefore, the premise of the question is not necessarily true.'
$ws.Range("F6").Value = '{''Explanation:'': ''The question assumes that if a problem can be quickly verified, then it can also be quickly solved. This is the core of the P versus NP problem in computer science. However, it is not necessarily true that all quickly verifiable problems are also quickly solvable. In fact, many problems exist where the answer can be verified in polynomial time, but there is no known way to solve the original problem in polynomial time.'', ''Answer:'': ''Not all quickly verifiable problems are also quickly solvable.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://en.wikipedia.org/wiki/P_versus_NP_problem>'', ''Core Concept:'': ''P versus NP problem'', ''Premise of the Question:'': ''Invalid. The premise of the question assumes that quickly verifiable problems are also quickly solvable, but this is not necessarily true.''}'
$ws.Range("G6").Value = '{''Explanation:'': ''The premise of the question assumes that if a problem can be quickly verified, then it can also be quickly solved. However, this is the fundamental question that the P versus NP problem tries to answer. In the field of computer science, a problem is considered to be in NP if a potential solution can be verified quickly, but not necessarily solved quickly. Therefore, the premise of the question is not necessarily true.'', ''Answer:'': ''Not all quickly verifiable problems are also quickly solvable.'', ''Confidence Level:'': ''100.0%'', ''Source:'': ''P versus NP problem (<https://en.wikipedia.org/wiki/P_versus_NP_problem>)'', ''Core Concept:'': ''NP-completeness'', ''Premise of the Question:'': ''Invalid. The premise assumes that quickly verifiable problems can also be quickly solved, but this is not always the case, as demonstrated by the P versus NP problem.'', ''Certainty_Estimation'': ''Yes, both responses are discussing the same topic. They both address the idea that not all problems which can be quickly verified can also be quickly solved, and refer to the P versus NP problem in computer science to support their argument.''}'
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = '{''QueryID'': ''4'', ''related_questions'': [{''name'': ''Are all problems solvable in polynomial time?'', ''url'': ''https://www.cs.umd.edu/class/fall2017/cmsc451-0101/Handouts/sol5.pdf'', ''snippet'': ''(i) All NP-complete problems are solvable in polynomial time: Yes. Every problem in NP is polynomially reducible to SAT, and SAT is reducible to every NP-hard problem. Therefore, a polynomial time solution to any NP-hard problem (such as 3Col) implies that every problem in NP can be solved in polynomial time.''}, {''name'': ''Why can we solve every problem in NP in exponential time?'', ''url'': ''https://www.quora.com/Can-all-NP-problems-be-solved-in-exponential-time#:~:text=Yes%2C%20cause%20the%20class%20of,in%20that%20computational%20complexity%20class).'', ''snippet'': ''Yes, cause the class of NP is contained within the class EXP (EXPTIME), while EXP is the class of all problems solvable (decidable) within exponential running-time (relative to the representation of the size of the non-unary input to the problem in that computational complexity class).''}, {''name'': ''What does it mean to be solvable in polynomial time?'', ''url'': ''https://www.britannica.com/science/polynomial-time-algorithm#:~:text=computational%20problems&text=%E2%80%A6can%20be%20solved%20in%20%E2%80%9Cpolynomial,the%20input%20for%20the%20problem.'', ''snippet'': ''computational problems …can be solved in “polynomial time,” which means that an algorithm exists for its solution such that the number of steps in the algorithm is bounded by a polynomial function of n, where n corresponds to the length of the input for the problem.''}, {''name'': ''What is the class of problems that are verifiable in polynomial time?'', ''url'': ''https://en.wikipedia.org/wiki/P_versus_NP_problem#:~:text=The%20class%20of%20questions%20where,for%20%22nondeterministic%20polynomial%20time%22.'', ''snippet'': ''The class of questions where an answer can be verified in polynomial time is NP, standing for "nondeterministic polynomial time".''}], ''organic_results'': [{''name'': ''P versus NP problem'', ''url'': ''https://en.wikipedia.org/wiki/P_versus_NP_problem'', ''snippet'': ''Informally, it asks whether every problem whose solution can be quickly verified can also be quickly solved. Here, quickly means an algorithm that solves ...''}, {''name'': ''Is there a task that is solvable in polynomial time but not ...'', ''url'': ''https://www.quora.com/Is-there-a-task-that-is-solvable-in-polynomial-time-but-not-verifiable-in-polynomial-time'', ''snippet'': ''For some problems, the answer can be verified to be correct in Polynomial Time, even if there is no known way of solving the original problem in ...''}]}'
$ws.Range("J6").Value = '{''Explanation:'': ''The question assumes that if a problem can be quickly verified, then it can also be quickly solved. This is the central question of the P versus NP problem in computer science, where P represents the problems that can be quickly solved, and NP represents the problems that can be quickly verified. The question is whether these two classes of problems are equivalent.'', ''Answer:'': ''It is not necessarily true that all problems that can be quickly verified can also be quickly solved.'', ''Confidence Level:'': ''90%'', ''Source:'': ''P versus NP problem (<https://en.wikipedia.org/wiki/P_versus_NP_problem>)'', ''Core Concept:'': ''P versus NP problem'', ''Premise of the Question:'': ''Invalid. The premise of the question is the fundamental question being asked in the P versus NP problem, and it is not necessarily true.'', ''Certainty_Estimation'': ''Yes, both responses are discussing the same topic. They both address the P versus NP problem in computer science, which involves the relationship between problems that can be quickly solved and those that can be quickly verified. The main subject, context, and key points in both responses are essentially the same.''}'
$ws.Range("K6").Value = '{''Explanation:'': ''The claim that all problems which are quickly verifiable are also quickly solvable is related to the P versus NP problem in computer science, which concerns whether every problem whose solution can be quickly verified can also be quickly solved. Here, quickly means in polynomial time, i.e., the time it takes to solve the problem increases at most quadratically with the size of the input. However, it is currently unknown whether P equals NP, and many computer scientists believe that they are not equal, implying that there are problems that can be quickly verified but not quickly solved.'', ''Answer:'': ''It is not necessarily true that all problems which are quickly verifiable are also quickly solvable.'', ''Confidence Level:'': ''90%'', ''Source:'': ''<https://en.wikipedia.org/wiki/P_versus_NP_problem>'', ''Core Concept:'': ''P versus NP problem'', ''Premise of the Question:'': ''Invalid. The premise assumes that quickly verifiable problems are also quickly solvable, but this is currently unknown and many experts believe it to be false.'', ''Certainty_Estimation'': ''Yes.\n\nConfidence: 100%''}'
$ws.Range("L6").Value = '{''Explanation:'': ''The premise of the question assumes that if a problem can be quickly verified, then it can also be quickly solved. However, this is the fundamental question that the P versus NP problem tries to answer. In the field of computer science, a problem is considered to be in NP if a potential solution can be verified quickly, but not necessarily solved quickly. Therefore, the premise of the question is not necessarily true.'', ''Answer:'': ''Not all quickly verifiable problems are also quickly solvable.'', ''Confidence Level:'': ''90%'', ''Source:'': ''P versus NP problem (<https://en.wikipedia.org/wiki/P_versus_NP_problem>)'', ''Core Concept:'': ''NP-completeness'', ''Premise of the Question:'': ''Invalid. The premise assumes that quickly verifiable problems can also be quickly solved, but this is not always the case, as demonstrated by the P versus NP problem.'', ''Certainty_Estimation'': ''Yes, both responses are discussing the same topic. They both address the idea that not all problems which can be quickly verified can also be quickly solved, and refer to the P versus NP problem in computer science to support their argument.''}'
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("O6").Value = 'Correct'
$ws.Range("P6").Value = 'The response accurately addresses the misconception in the question and correctly identifies the relationship between quickly verifiable and quickly solvable problems as the P versus NP problem, which remains an open question in computer science. Thus, the primary answer provided by the response is accurate.'
